$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename header cells on existing sheets
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Data rows 2-55 (54 rows)
$data = New-Object "object[,]" 54,4
$data[0,0] = 44941.99999999999
$data[0,1] = 4
$data[0,2] = -5.241064555832582
$data[0,3] = 12.88500144853277
$data[1,0] = 44948.99999999999
$data[1,1] = 4
$data[1,2] = -5.279303743595799
$data[1,3] = 13.20605555173696
$data[2,0] = 44955.99999999999
$data[2,1] = 4
$data[2,2] = -4.738573418130684
$data[2,3] = 13.303374038091
$data[3,0] = 44962.99999999999
$data[3,1] = 4
$data[3,2] = -5.159647122564289
$data[3,3] = 13.55400100876084
$data[4,0] = 44969.99999999999
$data[4,1] = 4
$data[4,2] = -5.122983693816792
$data[4,3] = 14.04930067289999
$data[5,0] = 44976.99999999999
$data[5,1] = 4
$data[5,2] = -4.255297885698592
$data[5,3] = 13.34098135130529
$data[6,0] = 44990.99999999999
$data[6,1] = 5
$data[6,2] = -4.150574472275073
$data[6,3] = 14.31932792743271
$data[7,0] = 44997.99999999999
$data[7,1] = 5
$data[7,2] = -4.451392782584371
$data[7,3] = 13.70894123319413
$data[8,0] = 45011.99999999999
$data[8,1] = 5
$data[8,2] = -4.321804542456886
$data[8,3] = 13.50444843221717
$data[9,0] = 45025.99999999999
$data[9,1] = 5
$data[9,2] = -3.710504407426772
$data[9,3] = 13.59463288880438
$data[10,0] = 45046.99999999999
$data[10,1] = 5
$data[10,2] = -3.336997928099318
$data[10,3] = 14.10793109772424
$data[11,0] = 45053.99999999999
$data[11,1] = 5
$data[11,2] = -3.224830691526629
$data[11,3] = 14.10610967667534
$data[12,0] = 45060.99999999999
$data[12,1] = 5
$data[12,2] = -3.885422268152837
$data[12,3] = 13.86973087041261
$data[13,0] = 45067.99999999999
$data[13,1] = 5
$data[13,2] = -3.88379982421427
$data[13,3] = 14.07965793592481
$data[14,0] = 45074.99999999999
$data[14,1] = 5
$data[14,2] = -3.830321411224957
$data[14,3] = 14.51266426752773
$data[15,0] = 45081.99999999999
$data[15,1] = 6
$data[15,2] = -2.619024426021853
$data[15,3] = 14.43314591752319
$data[16,0] = 45088.99999999999
$data[16,1] = 6
$data[16,2] = -3.133095809338097
$data[16,3] = 15.41617399092681
$data[17,0] = 45116.99999999999
$data[17,1] = 6
$data[17,2] = -3.609069070069317
$data[17,3] = 15.0262526857956
$data[18,0] = 45123.99999999999
$data[18,1] = 6
$data[18,2] = -3.392408365937464
$data[18,3] = 14.09819702159873
$data[19,0] = 45130.99999999999
$data[19,1] = 6
$data[19,2] = -3.202228025960757
$data[19,3] = 15.67122722211109
$data[20,0] = 45158.99999999999
$data[20,1] = 6
$data[20,2] = -2.643060311282154
$data[20,3] = 15.95074155976132
$data[21,0] = 45165.99999999999
$data[21,1] = 6
$data[21,2] = -2.288768784886306
$data[21,3] = 15.42441193804697
$data[22,0] = 45179.99999999999
$data[22,1] = 7
$data[22,2] = -3.047153937250428
$data[22,3] = 16.41168292617332
$data[23,0] = 45186.99999999999
$data[23,1] = 7
$data[23,2] = -2.970877608680877
$data[23,3] = 16.18151385377156
$data[24,0] = 45193.99999999999
$data[24,1] = 7
$data[24,2] = -2.275565042524093
$data[24,3] = 16.08033076885601
$data[25,0] = 45214.99999999999
$data[25,1] = 7
$data[25,2] = -1.878651328483616
$data[25,3] = 15.89435564301814
$data[26,0] = 45221.99999999999
$data[26,1] = 7
$data[26,2] = -2.395637758983164
$data[26,3] = 16.76970727080994
$data[27,0] = 45228.99999999999
$data[27,1] = 7
$data[27,2] = -2.146359271350792
$data[27,3] = 16.25034105480453
$data[28,0] = 45242.99999999999
$data[28,1] = 7
$data[28,2] = -1.368748534895342
$data[28,3] = 16.38778308336986
$data[29,0] = 45249.99999999999
$data[29,1] = 7
$data[29,2] = -1.97964303843687
$data[29,3] = 16.37396626206261
$data[30,0] = 45256.99999999999
$data[30,1] = 7
$data[30,2] = -1.911993395322696
$data[30,3] = 16.16707095220363
$data[31,0] = 45263.99999999999
$data[31,1] = 7
$data[31,2] = -1.565293897200069
$data[31,3] = 16.65298346869167
$data[32,0] = 45270.99999999999
$data[32,1] = 7
$data[32,2] = -1.286734103755333
$data[32,3] = 16.48732318981288
$data[33,0] = 45277.99999999999
$data[33,1] = 8
$data[33,2] = -1.904913348094943
$data[33,3] = 16.31285587753036
$data[34,0] = 45298.99999999999
$data[34,1] = 8
$data[34,2] = -0.7733744477293844
$data[34,3] = 16.97213812942804
$data[35,0] = 45396.99999999999
$data[35,1] = 9
$data[35,2] = -0.5399016817845228
$data[35,3] = 17.99919042609522
$data[36,0] = 45410.99999999999
$data[36,1] = 9
$data[36,2] = 0.06907701779345735
$data[36,3] = 17.69433011342685
$data[37,0] = 45417.99999999999
$data[37,1] = 9
$data[37,2] = -0.08804307294050424
$data[37,3] = 18.25963753587872
$data[38,0] = 45424.99999999999
$data[38,1] = 9
$data[38,2] = -0.5684372820410729
$data[38,3] = 18.81868522856785
$data[39,0] = 45431.99999999999
$data[39,1] = 9
$data[39,2] = 0.09377849835529364
$data[39,3] = 18.04596019668685
$data[40,0] = 45438.99999999999
$data[40,1] = 9
$data[40,2] = 0.1086722882481971
$data[40,3] = 18.50621890675455
$data[41,0] = 45445.99999999999
$data[41,1] = 9
$data[41,2] = 0.6950842146994221
$data[41,3] = 17.97808475585292
$data[42,0] = 45452.99999999999
$data[42,1] = 9
$data[42,2] = 0.7042749185506225
$data[42,3] = 17.76367779022069
$data[43,0] = 45459.99999999999
$data[43,1] = 9
$data[43,2] = 0.3130393590887959
$data[43,3] = 18.52693475470387
$data[44,0] = 45487.99999999999
$data[44,1] = 10
$data[44,2] = -0.2698686961861273
$data[44,3] = 18.17315823492288
$data[45,0] = 45543.99999999999
$data[45,1] = 10
$data[45,2] = 1.25520199098237
$data[45,3] = 18.35310099401703
$data[46,0] = 45550.99999999999
$data[46,1] = 10
$data[46,2] = 1.287205102491998
$data[46,3] = 19.13645774855748
$data[47,0] = 45557.99999999999
$data[47,1] = 10
$data[47,2] = 0.8926048962780404
$data[47,3] = 19.14573988672929
$data[48,0] = 45564.99999999999
$data[48,1] = 10
$data[48,2] = 1.630386498242182
$data[48,3] = 19.29853492949341
$data[49,0] = 45571.99999999999
$data[49,1] = 11
$data[49,2] = 0.9343135652180586
$data[49,3] = 19.66341129676436
$data[50,0] = 45578.99999999999
$data[50,1] = 11
$data[50,2] = 1.215206612480937
$data[50,3] = 19.81068059492162
$data[51,0] = 45585.99999999999
$data[51,1] = 11
$data[51,2] = 1.922381761388726
$data[51,3] = 20.23711162706165
$data[52,0] = 45592.99999999999
$data[52,1] = 11
$data[52,2] = 1.508602675324569
$data[52,3] = 19.77289187258874
$data[53,0] = 45599.99999999999
$data[53,1] = 11
$data[53,2] = 2.098452748284541
$data[53,3] = 19.74274621743338

$ws3.Range("A2:D55").Value = $data

# Copy header + date formatting from the existing "Weekly Quantity" sheet
$ws1.Range("A1:B1").Copy() | Out-Null
$ws3.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Copy() | Out-Null
$ws3.Range("A2:A55").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
